# Update cryptocurrency price (D) and 1h volume change (E) figures to the
# latest values pulled by the GitHub Actions symbol-list refresh job.
# Each numeric-looking value is entered as literal text (matching the
# workbook's existing inline-string cell type for these columns), so we
# force the Text number format before assigning, then restore the cell
# style back to Normal so only the displayed value changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "331.92"
Set-TextValue "E2" "0.27%"
Set-TextValue "D3" "41.69"
Set-TextValue "E3" "3.88%"
Set-TextValue "D4" "5.658"
Set-TextValue "E4" "-1.07%"
Set-TextValue "D5" "0.08361"
Set-TextValue "E5" "3.17%"
Set-TextValue "D6" "8.790"
Set-TextValue "E6" "1.59%"
Set-TextValue "D7" "2.012"
Set-TextValue "E7" "2.73%"
Set-TextValue "D8" "4.550"
Set-TextValue "E8" "1.45%"
Set-TextValue "E9" "1.21%"
Set-TextValue "D10" "0.9262"
Set-TextValue "E10" "-0.33%"
Set-TextValue "D11" "0.1288"
Set-TextValue "E11" "2.97%"
Set-TextValue "D12" "0.1960"
Set-TextValue "E12" "-0.09%"
Set-TextValue "D13" "0.09406"
Set-TextValue "E13" "1.92%"
Set-TextValue "D14" "0.03884"
Set-TextValue "E14" "4.11%"
Set-TextValue "E15" "0.95%"
Set-TextValue "D16" "0.001312"
Set-TextValue "E16" "1.49%"
Set-TextValue "D17" "0.006113"
Set-TextValue "E17" "-3.20%"
Set-TextValue "D18" "3.438"
Set-TextValue "E18" "1.84%"
Set-TextValue "E19" "1.32%"
Set-TextValue "D20" "7.978"
Set-TextValue "E20" "-9.00%"
Set-TextValue "D21" "0.1373"
Set-TextValue "E21" "0.53%"
Set-TextValue "D22" "0.2612"
Set-TextValue "E22" "0.48%"
Set-TextValue "D23" "0.04421"
Set-TextValue "E23" "0.09%"
Set-TextValue "D24" "0.001254"
Set-TextValue "E24" "0.48%"
Set-TextValue "D25" "0.004463"
Set-TextValue "E25" "0.27%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "-2.86%"
Set-TextValue "D39" "0.02795"
Set-TextValue "E39" "1.87%"
Set-TextValue "D40" "0.05565"
Set-TextValue "E40" "0.23%"
Set-TextValue "D41" "0.007794"
Set-TextValue "E41" "3.78%"
Set-TextValue "D42" "0.1434"
Set-TextValue "E42" "0.86%"
Set-TextValue "D43" "0.009319"
Set-TextValue "E43" "-5.02%"
Set-TextValue "D44" "0.002157"
Set-TextValue "E44" "2.79%"
Set-TextValue "E45" "-6.40%"
Set-TextValue "D46" "0.00007012"
Set-TextValue "E46" "3.72%"
Set-TextValue "E47" "0.38%"
Set-TextValue "D48" "0.003499"
Set-TextValue "E48" "14.29%"
Set-TextValue "D49" "0.002281"
Set-TextValue "E49" "0.49%"
Set-TextValue "E50" "0.38%"
Set-TextValue "E51" "0.38%"
